# PUB-2917 Add Planning Court in London Administrative Court
#
# The workbook originally has a single sheet ("Sheet1") holding the London
# administrative court cause-list data. This change duplicates that sheet so
# the workbook ends up with two tabs sharing the same layout/data:
#   1. "London administrative court" - a brand new copy of the data, placed
#      before the original sheet and left as the active/selected tab.
#   2. "Planning court" - the original sheet, renamed, no longer selected.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item(1)

# Copy the sheet, inserting the duplicate before the source sheet (mirrors
# Excel's "Move or Copy... > Create a copy" placed ahead of the original).
$sourceSheet.Copy($sourceSheet)

# Re-fetch fresh references by position now that the workbook has two sheets -
# the new copy takes tab 1 and becomes the active sheet, the original slides
# down to tab 2.
$newSheet = $wb.Worksheets.Item(1)
$originalSheet = $wb.Worksheets.Item(2)

$newSheet.Name = "London administrative court"
$originalSheet.Name = "Planning court"

# Match the author's recorded selection on the new, active tab.
$newSheet.Activate()
$newSheet.Range("D35").Select() | Out-Null
